# Add a new "FTHC_Average_Demand5" results sheet at the end of the workbook,
# mirroring the layout of the existing FTNC_Average_Demand5* sheets
# (header row with In-vehicle/At-stop/Extra/Tardiness/Total, one data row).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet so it lands at the
# end of the tab order (Worksheets.Add() with no args inserts at the front).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "FTHC_Average_Demand5"

# Copy the header row + row formatting/styles from an existing sibling sheet
# so the new sheet matches the established look (bold, centered, bordered
# cells) instead of re-deriving a brand-new style.
$srcWs = $wb.Worksheets.Item("FTNC_Average_Demand54")
$srcWs.Range("B1:F1").Copy($ws.Range("B1"))
$srcWs.Range("A2:F2").Copy($ws.Range("A2"))

# Overwrite row 2 with this sheet's own label + computed averages.
$ws.Range("A2").Value = "FTHC_Average_Demand_5"
$ws.Range("B2").Value = 12.11441948061514
$ws.Range("C2").Value = 182.7178957992942
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 7.255286811132406
$ws.Range("F2").Value = 202.0876020910418
